$wb = $excel.ActiveWorkbook

$edits = @(
    @{Sheet="ALC"; Cell="H32"; Op="Set"; Value=3274},
    @{Sheet="ALC"; Cell="I32"; Op="Set"; Value=2666.6667},
    @{Sheet="ALC"; Cell="K32"; Op="Set"; Value=2666.6667},
    @{Sheet="ALC"; Cell="M32"; Op="Set"; Value=-2340.6667},
    @{Sheet="ALC"; Cell="H34"; Op="Set"; Value=10000},
    @{Sheet="ALC"; Cell="I34"; Op="Set"; Value=10000},
    @{Sheet="ALC"; Cell="K34"; Op="Set"; Value=10000},
    @{Sheet="ALC"; Cell="M34"; Op="Set"; Value=-9797},
    @{Sheet="ALC"; Cell="H36"; Op="Set"; Value=10000},
    @{Sheet="ALC"; Cell="I36"; Op="Set"; Value=10000},
    @{Sheet="ALC"; Cell="K36"; Op="Set"; Value=10000},
    @{Sheet="ALC"; Cell="M36"; Op="Set"; Value=-9285},
    @{Sheet="ALC"; Cell="H41"; Op="Set"; Value=997},
    @{Sheet="ALC"; Cell="J41"; Op="Set"; Value=0},
    @{Sheet="ALC"; Cell="L41"; Op="Set"; Value=0},
    @{Sheet="ALC"; Cell="N41"; Op="Clear"},
    @{Sheet="ALC"; Cell="H43"; Op="Set"; Value=4285.4116},
    @{Sheet="ALC"; Cell="I43"; Op="Set"; Value=1279.8},
    @{Sheet="ALC"; Cell="J43"; Op="Set"; Value=5537.75},
    @{Sheet="ALC"; Cell="K43"; Op="Set"; Value=1279.8},
    @{Sheet="ALC"; Cell="L43"; Op="Set"; Value=5537.75},
    @{Sheet="ALC"; Cell="M43"; Op="Set"; Value=-1210.8},
    @{Sheet="ALC"; Cell="N43"; Op="Set"; Value=-5675.75},
    @{Sheet="ALC"; Cell="H64"; Op="Set"; Value=5204.197},
    @{Sheet="ALC"; Cell="I64"; Op="Set"; Value=3261.5945},
    @{Sheet="ALC"; Cell="J64"; Op="Set"; Value=7682.6895},
    @{Sheet="ALC"; Cell="K64"; Op="Set"; Value=3261.5945},
    @{Sheet="ALC"; Cell="L64"; Op="Set"; Value=7682.6895},
    @{Sheet="ALC"; Cell="M64"; Op="Set"; Value=-3013.5945},
    @{Sheet="ALC"; Cell="N64"; Op="Set"; Value=-8178.6895},
    @{Sheet="ALC"; Cell="H67"; Op="Set"; Value=5204.197},
    @{Sheet="ALC"; Cell="I67"; Op="Set"; Value=3261.5945},
    @{Sheet="ALC"; Cell="J67"; Op="Set"; Value=7682.6895},
    @{Sheet="ALC"; Cell="K67"; Op="Set"; Value=3261.5945},
    @{Sheet="ALC"; Cell="L67"; Op="Set"; Value=7682.6895},
    @{Sheet="ALC"; Cell="M67"; Op="Set"; Value=-2403.5945},
    @{Sheet="ALC"; Cell="N67"; Op="Set"; Value=-9398.6895},
    @{Sheet="ALC"; Cell="H88"; Op="Set"; Value=3889.2144},
    @{Sheet="ALC"; Cell="I88"; Op="Set"; Value=551.5},
    @{Sheet="ALC"; Cell="K88"; Op="Set"; Value=551.5},
    @{Sheet="ALC"; Cell="M88"; Op="Set"; Value=-145.5},
    @{Sheet="ALC"; Cell="H91"; Op="Set"; Value=3889.2144},
    @{Sheet="ALC"; Cell="I91"; Op="Set"; Value=551.5},
    @{Sheet="ALC"; Cell="K91"; Op="Set"; Value=551.5},
    @{Sheet="ALC"; Cell="M91"; Op="Set"; Value=852.5},
    @{Sheet="ALC"; Cell="H100"; Op="Set"; Value=6741.3335},
    @{Sheet="ALC"; Cell="I100"; Op="Set"; Value=3750},
    @{Sheet="ALC"; Cell="K100"; Op="Set"; Value=3750},
    @{Sheet="ALC"; Cell="M100"; Op="Set"; Value=-3209},
    @{Sheet="ALC"; Cell="H116"; Op="Set"; Value=4316.6665},
    @{Sheet="ALC"; Cell="I116"; Op="Set"; Value=3395},
    @{Sheet="ALC"; Cell="J116"; Op="Set"; Value=5054},
    @{Sheet="ALC"; Cell="K116"; Op="Set"; Value=3395},
    @{Sheet="ALC"; Cell="L116"; Op="Set"; Value=5054},
    @{Sheet="ALC"; Cell="M116"; Op="Set"; Value=47},
    @{Sheet="ALC"; Cell="N116"; Op="Set"; Value=-11938},
    @{Sheet="ALC"; Cell="H118"; Op="Set"; Value=2073.8},
    @{Sheet="ALC"; Cell="I118"; Op="Set"; Value=1193.1111},
    @{Sheet="ALC"; Cell="K118"; Op="Set"; Value=3579.3333},
    @{Sheet="ALC"; Cell="M118"; Op="Set"; Value=-1922.3333},
    @{Sheet="ALC"; Cell="H132"; Op="Set"; Value=1660.2727},
    @{Sheet="ALC"; Cell="I132"; Op="Set"; Value=1401.3},
    @{Sheet="ALC"; Cell="K132"; Op="Set"; Value=4203.9},
    @{Sheet="ALC"; Cell="M132"; Op="Set"; Value=-1673.9},
    @{Sheet="ALC"; Cell="H141"; Op="Set"; Value=8687.053},
    @{Sheet="ALC"; Cell="I141"; Op="Set"; Value=10766.25},
    @{Sheet="ALC"; Cell="J141"; Op="Set"; Value=8132.6},
    @{Sheet="ALC"; Cell="K141"; Op="Set"; Value=32298.75},
    @{Sheet="ALC"; Cell="L141"; Op="Set"; Value=24397.8},
    @{Sheet="ALC"; Cell="M141"; Op="Set"; Value=-27118.75},
    @{Sheet="ALC"; Cell="N141"; Op="Set"; Value=-34757.8},
    @{Sheet="ARM"; Cell="H32"; Op="Set"; Value=3570.394},
    @{Sheet="ARM"; Cell="I32"; Op="Set"; Value=3672.5625},
    @{Sheet="ARM"; Cell="K32"; Op="Set"; Value=3672.5625},
    @{Sheet="ARM"; Cell="M32"; Op="Set"; Value=-3385.5625},
    @{Sheet="ARM"; Cell="H45"; Op="Set"; Value=1629.5},
    @{Sheet="ARM"; Cell="I45"; Op="Set"; Value=1629.5},
    @{Sheet="ARM"; Cell="J45"; Op="Set"; Value=0},
    @{Sheet="ARM"; Cell="K45"; Op="Set"; Value=1629.5},
    @{Sheet="ARM"; Cell="L45"; Op="Set"; Value=0},
    @{Sheet="ARM"; Cell="M45"; Op="Set"; Value=-1252.5},
    @{Sheet="ARM"; Cell="N45"; Op="Clear"},
    @{Sheet="ARM"; Cell="H61"; Op="Set"; Value=5684.2354},
    @{Sheet="ARM"; Cell="I61"; Op="Set"; Value=6434.231},
    @{Sheet="ARM"; Cell="J61"; Op="Set"; Value=3246.75},
    @{Sheet="ARM"; Cell="K61"; Op="Set"; Value=6434.231},
    @{Sheet="ARM"; Cell="L61"; Op="Set"; Value=3246.75},
    @{Sheet="ARM"; Cell="M61"; Op="Set"; Value=-6222.231},
    @{Sheet="ARM"; Cell="N61"; Op="Set"; Value=-3670.75},
    @{Sheet="ARM"; Cell="H74"; Op="Set"; Value=2501.9443},
    @{Sheet="ARM"; Cell="I74"; Op="Set"; Value=2659.25},
    @{Sheet="ARM"; Cell="K74"; Op="Set"; Value=2659.25},
    @{Sheet="ARM"; Cell="M74"; Op="Set"; Value=-1785.25},
    @{Sheet="ARM"; Cell="H77"; Op="Set"; Value=2501.9443},
    @{Sheet="ARM"; Cell="I77"; Op="Set"; Value=2659.25},
    @{Sheet="ARM"; Cell="K77"; Op="Set"; Value=13296.25},
    @{Sheet="ARM"; Cell="M77"; Op="Set"; Value=-8928.25},
    @{Sheet="ARM"; Cell="H122"; Op="Set"; Value=1532.8636},
    @{Sheet="ARM"; Cell="I122"; Op="Set"; Value=1250.4857},
    @{Sheet="ARM"; Cell="K122"; Op="Set"; Value=3751.4571},
    @{Sheet="ARM"; Cell="M122"; Op="Set"; Value=-1301.4571},
    @{Sheet="ARM"; Cell="H136"; Op="Set"; Value=5684.2354},
    @{Sheet="ARM"; Cell="I136"; Op="Set"; Value=6434.231},
    @{Sheet="ARM"; Cell="J136"; Op="Set"; Value=3246.75},
    @{Sheet="ARM"; Cell="K136"; Op="Set"; Value=19302.693},
    @{Sheet="ARM"; Cell="L136"; Op="Set"; Value=9740.25},
    @{Sheet="ARM"; Cell="M136"; Op="Set"; Value=-16752.693},
    @{Sheet="ARM"; Cell="N136"; Op="Set"; Value=-14840.25},
    @{Sheet="BSM"; Cell="H134"; Op="Set"; Value=2764.077},
    @{Sheet="BSM"; Cell="I134"; Op="Set"; Value=2721.182},
    @{Sheet="BSM"; Cell="K134"; Op="Set"; Value=8163.545999999999},
    @{Sheet="BSM"; Cell="M134"; Op="Set"; Value=-5628.545999999999},
    @{Sheet="CRP"; Cell="H134"; Op="Set"; Value=3001.1177},
    @{Sheet="CRP"; Cell="I134"; Op="Set"; Value=3476.125},
    @{Sheet="CRP"; Cell="J134"; Op="Set"; Value=2578.889},
    @{Sheet="CRP"; Cell="K134"; Op="Set"; Value=10428.375},
    @{Sheet="CRP"; Cell="L134"; Op="Set"; Value=7736.667},
    @{Sheet="CRP"; Cell="M134"; Op="Set"; Value=-7893.375},
    @{Sheet="CRP"; Cell="N134"; Op="Set"; Value=-12806.667},
    @{Sheet="CUL"; Cell="H38"; Op="Set"; Value=676.4737},
    @{Sheet="CUL"; Cell="J38"; Op="Set"; Value=976.38464},
    @{Sheet="CUL"; Cell="L38"; Op="Set"; Value=2929.15392},
    @{Sheet="CUL"; Cell="N38"; Op="Set"; Value=-3623.15392},
    @{Sheet="CUL"; Cell="H81"; Op="Set"; Value=6544.1665},
    @{Sheet="CUL"; Cell="I81"; Op="Set"; Value=632.5},
    @{Sheet="CUL"; Cell="K81"; Op="Set"; Value=1897.5},
    @{Sheet="CUL"; Cell="M81"; Op="Set"; Value=-774.5},
    @{Sheet="CUL"; Cell="H84"; Op="Set"; Value=6544.1665},
    @{Sheet="CUL"; Cell="I84"; Op="Set"; Value=632.5},
    @{Sheet="CUL"; Cell="K84"; Op="Set"; Value=5692.5},
    @{Sheet="CUL"; Cell="M84"; Op="Set"; Value=-76.5},
    @{Sheet="CUL"; Cell="H100"; Op="Set"; Value=66700},
    @{Sheet="CUL"; Cell="I100"; Op="Set"; Value=0},
    @{Sheet="CUL"; Cell="K100"; Op="Set"; Value=0},
    @{Sheet="CUL"; Cell="M100"; Op="Clear"},
    @{Sheet="CUL"; Cell="H112"; Op="Set"; Value=146123.86},
    @{Sheet="CUL"; Cell="I112"; Op="Set"; Value=250541.5},
    @{Sheet="CUL"; Cell="J112"; Op="Set"; Value=6900.3335},
    @{Sheet="CUL"; Cell="K112"; Op="Set"; Value=751624.5},
    @{Sheet="CUL"; Cell="L112"; Op="Set"; Value=20701.0005},
    @{Sheet="CUL"; Cell="M112"; Op="Set"; Value=-750516.5},
    @{Sheet="CUL"; Cell="N112"; Op="Set"; Value=-22917.0005},
    @{Sheet="CUL"; Cell="H113"; Op="Set"; Value=1216.0625},
    @{Sheet="CUL"; Cell="J113"; Op="Set"; Value=1331.3077},
    @{Sheet="CUL"; Cell="L113"; Op="Set"; Value=3993.9231},
    @{Sheet="CUL"; Cell="N113"; Op="Set"; Value=-8333.9231},
    @{Sheet="GSM"; Cell="H36"; Op="Set"; Value=6008.25},
    @{Sheet="GSM"; Cell="I36"; Op="Set"; Value=6008.25},
    @{Sheet="GSM"; Cell="J36"; Op="Set"; Value=0},
    @{Sheet="GSM"; Cell="K36"; Op="Set"; Value=6008.25},
    @{Sheet="GSM"; Cell="L36"; Op="Set"; Value=0},
    @{Sheet="GSM"; Cell="M36"; Op="Set"; Value=-5523.25},
    @{Sheet="GSM"; Cell="N36"; Op="Clear"},
    @{Sheet="GSM"; Cell="H43"; Op="Set"; Value=9000},
    @{Sheet="GSM"; Cell="I43"; Op="Set"; Value=9000},
    @{Sheet="GSM"; Cell="K43"; Op="Set"; Value=9000},
    @{Sheet="GSM"; Cell="M43"; Op="Set"; Value=-8849},
    @{Sheet="GSM"; Cell="H102"; Op="Set"; Value=1379.8182},
    @{Sheet="GSM"; Cell="I102"; Op="Set"; Value=1402.762},
    @{Sheet="GSM"; Cell="K102"; Op="Set"; Value=1402.762},
    @{Sheet="GSM"; Cell="M102"; Op="Set"; Value=219.2380000000001},
    @{Sheet="GSM"; Cell="H122"; Op="Set"; Value=5326.6562},
    @{Sheet="GSM"; Cell="I122"; Op="Set"; Value=5194.8213},
    @{Sheet="GSM"; Cell="J122"; Op="Set"; Value=6249.5},
    @{Sheet="GSM"; Cell="K122"; Op="Set"; Value=15584.4639},
    @{Sheet="GSM"; Cell="L122"; Op="Set"; Value=18748.5},
    @{Sheet="GSM"; Cell="M122"; Op="Set"; Value=-13134.4639},
    @{Sheet="GSM"; Cell="N122"; Op="Set"; Value=-23648.5},
    @{Sheet="GSM"; Cell="H126"; Op="Set"; Value=3499.5454},
    @{Sheet="GSM"; Cell="I126"; Op="Set"; Value=3152.6},
    @{Sheet="GSM"; Cell="K126"; Op="Set"; Value=9457.799999999999},
    @{Sheet="GSM"; Cell="M126"; Op="Set"; Value=-6987.799999999999},
    @{Sheet="GSM"; Cell="H132"; Op="Set"; Value=2374},
    @{Sheet="GSM"; Cell="I132"; Op="Set"; Value=2588.7},
    @{Sheet="GSM"; Cell="J132"; Op="Set"; Value=1837.25},
    @{Sheet="GSM"; Cell="K132"; Op="Set"; Value=7766.099999999999},
    @{Sheet="GSM"; Cell="L132"; Op="Set"; Value=5511.75},
    @{Sheet="GSM"; Cell="M132"; Op="Set"; Value=-5236.099999999999},
    @{Sheet="GSM"; Cell="N132"; Op="Set"; Value=-10571.75},
    @{Sheet="LTW"; Cell="H16"; Op="Set"; Value=654.2857},
    @{Sheet="LTW"; Cell="J16"; Op="Set"; Value=433.33334},
    @{Sheet="LTW"; Cell="L16"; Op="Set"; Value=433.33334},
    @{Sheet="LTW"; Cell="N16"; Op="Set"; Value=-773.33334},
    @{Sheet="LTW"; Cell="H42"; Op="Set"; Value=35000},
    @{Sheet="LTW"; Cell="J42"; Op="Set"; Value=35000},
    @{Sheet="LTW"; Cell="L42"; Op="Set"; Value=35000},
    @{Sheet="LTW"; Cell="N42"; Op="Set"; Value=-36126},
    @{Sheet="LTW"; Cell="H46"; Op="Set"; Value=2582.8},
    @{Sheet="LTW"; Cell="I46"; Op="Set"; Value=918},
    @{Sheet="LTW"; Cell="K46"; Op="Set"; Value=918},
    @{Sheet="LTW"; Cell="M46"; Op="Set"; Value=-730},
    @{Sheet="LTW"; Cell="H49"; Op="Set"; Value=35000},
    @{Sheet="LTW"; Cell="J49"; Op="Set"; Value=35000},
    @{Sheet="LTW"; Cell="L49"; Op="Set"; Value=35000},
    @{Sheet="LTW"; Cell="N49"; Op="Set"; Value=-35294},
    @{Sheet="LTW"; Cell="H93"; Op="Set"; Value=3645.4075},
    @{Sheet="LTW"; Cell="I93"; Op="Set"; Value=1727.1666},
    @{Sheet="LTW"; Cell="J93"; Op="Set"; Value=5180},
    @{Sheet="LTW"; Cell="K93"; Op="Set"; Value=1727.1666},
    @{Sheet="LTW"; Cell="L93"; Op="Set"; Value=5180},
    @{Sheet="LTW"; Cell="M93"; Op="Set"; Value=-479.1666},
    @{Sheet="LTW"; Cell="N93"; Op="Set"; Value=-7676},
    @{Sheet="LTW"; Cell="H100"; Op="Set"; Value=4619.1333},
    @{Sheet="LTW"; Cell="I100"; Op="Set"; Value=2461},
    @{Sheet="LTW"; Cell="J100"; Op="Set"; Value=7085.5713},
    @{Sheet="LTW"; Cell="K100"; Op="Set"; Value=2461},
    @{Sheet="LTW"; Cell="L100"; Op="Set"; Value=7085.5713},
    @{Sheet="LTW"; Cell="M100"; Op="Set"; Value=-1920},
    @{Sheet="LTW"; Cell="N100"; Op="Set"; Value=-8167.5713},
    @{Sheet="LTW"; Cell="H122"; Op="Set"; Value=4224.696},
    @{Sheet="LTW"; Cell="I122"; Op="Set"; Value=4040.9285},
    @{Sheet="LTW"; Cell="J122"; Op="Set"; Value=4510.5557},
    @{Sheet="LTW"; Cell="K122"; Op="Set"; Value=12122.7855},
    @{Sheet="LTW"; Cell="L122"; Op="Set"; Value=13531.6671},
    @{Sheet="LTW"; Cell="M122"; Op="Set"; Value=-9672.7855},
    @{Sheet="LTW"; Cell="N122"; Op="Set"; Value=-18431.6671},
    @{Sheet="WVR"; Cell="H96"; Op="Set"; Value=3533.8333},
    @{Sheet="WVR"; Cell="I96"; Op="Set"; Value=1067.6666},
    @{Sheet="WVR"; Cell="K96"; Op="Set"; Value=1067.6666},
    @{Sheet="WVR"; Cell="M96"; Op="Set"; Value=305.3334},
    @{Sheet="WVR"; Cell="H113"; Op="Set"; Value=1856.8182},
    @{Sheet="WVR"; Cell="I113"; Op="Set"; Value=1747.5},
    @{Sheet="WVR"; Cell="J113"; Op="Set"; Value=2148.3333},
    @{Sheet="WVR"; Cell="K113"; Op="Set"; Value=5242.5},
    @{Sheet="WVR"; Cell="L113"; Op="Set"; Value=6444.999899999999},
    @{Sheet="WVR"; Cell="M113"; Op="Set"; Value=-3072.5},
    @{Sheet="WVR"; Cell="N113"; Op="Set"; Value=-10784.9999},
    @{Sheet="WVR"; Cell="H126"; Op="Set"; Value=1458},
    @{Sheet="WVR"; Cell="I126"; Op="Set"; Value=1112.1428},
    @{Sheet="WVR"; Cell="J126"; Op="Set"; Value=2265},
    @{Sheet="WVR"; Cell="K126"; Op="Set"; Value=3336.4284},
    @{Sheet="WVR"; Cell="L126"; Op="Set"; Value=6795},
    @{Sheet="WVR"; Cell="M126"; Op="Set"; Value=-866.4284000000002},
    @{Sheet="WVR"; Cell="N126"; Op="Set"; Value=-11735},
    @{Sheet="WVR"; Cell="H132"; Op="Set"; Value=4805.6816},
    @{Sheet="WVR"; Cell="I132"; Op="Set"; Value=4796.1904},
    @{Sheet="WVR"; Cell="K132"; Op="Set"; Value=14388.5712},
    @{Sheet="WVR"; Cell="M132"; Op="Set"; Value=-11858.5712},
)

foreach ($e in $edits) {
    $ws = $wb.Worksheets.Item($e.Sheet)
    if ($e.Op -eq "Clear") {
        $ws.Range($e.Cell).ClearContents()
    } else {
        $ws.Range($e.Cell).Value = $e.Value
    }
}

Write-Host "Applied $($edits.Count) cell edits (refreshed market price data)"
